$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-formatted cells to remain as text so Excel does not
# auto-convert plain numeric-looking strings (e.g. "0.524") into numbers.
$textCells = @("D5", "D6", "D9", "D11", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D35", "D36", "D37", "D39", "D40", "D41", "D43", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the source diff.
$ws.Range("D2").Value = "64.039.92"
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").Value = "3.140.67"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "611.40"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "146.38"
$ws.Range("E6").Value = "  -6.72%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.138.22"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("E10").Value = "  -6.72%  "
$ws.Range("D11").Value = "5.32"
$ws.Range("E11").Value = "  -7.93%  "
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("E13").Value = "  -7.27%  "
$ws.Range("D14").Value = "35.33"
$ws.Range("E14").Value = "  -9.35%  "
$ws.Range("D15").Value = "3.655.58"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "64.070.64"
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").Value = "3.139.39"
$ws.Range("E18").Value = "  -3.37%  "
$ws.Range("D19").Value = "6.86"
$ws.Range("E19").Value = "  -7.73%  "
$ws.Range("D20").Value = "476.13"
$ws.Range("E20").Value = "  -5.66%  "
$ws.Range("D21").Value = "14.71"
$ws.Range("E21").Value = "  -4.61%  "
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -6.25%  "
$ws.Range("D23").Value = "7.76"
$ws.Range("E23").Value = "  -3.99%  "
$ws.Range("D24").Value = "13.59"
$ws.Range("E24").Value = "  -7.39%  "
$ws.Range("D25").Value = "83.54"
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "2.81"
$ws.Range("E27").Value = "  -7.00%  "
$ws.Range("D28").Value = "8.41"
$ws.Range("E28").Value = "  -8.29%  "
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -8.55%  "
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("E31").Value = "  -10.95%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "2.72"
$ws.Range("E33").Value = "  -5.33%  "
$ws.Range("E34").Value = "  -6.21%  "
$ws.Range("D35").Value = "1.13"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -7.49%  "
$ws.Range("D37").Value = "53.46"
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("D38").Value = "0.0₃0732"
$ws.Range("E38").Value = "  -5.97%  "
$ws.Range("D39").Value = "461.15"
$ws.Range("E39").Value = "  -6.79%  "
$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -12.43%  "
$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  -6.74%  "
$ws.Range("E42").Value = "  -7.92%  "
$ws.Range("D43").Value = "8.40"
$ws.Range("E43").Value = "  -4.82%  "
$ws.Range("D44").Value = "2.845.15"
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("E45").Value = "  -9.03%  "
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  -10.40%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "26.31"
$ws.Range("E48").Value = "  -8.67%  "
$ws.Range("D49").Value = "2.36"
$ws.Range("E49").Value = "  -5.49%  "
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("D51").Value = "118.36"
$ws.Range("E51").Value = "  -2.04%  "
